# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" worksheets, reflecting newly scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new F value for "展览" sheet (sheet1)
$updates1 = @{
    3  = 1832
    5  = 795
    13 = 132
    14 = 150
    15 = 4291
    18 = 467
    19 = 409
    20 = 985
    21 = 1548
    23 = 43
    25 = 45
    26 = 2004
    27 = 67
    28 = 62
    29 = 2
    30 = 136
    31 = 56
    32 = 204
    33 = 27
}

# Row -> new F value for "全部类型" sheet (sheet4)
$updates4 = @{
    3  = 1832
    5  = 795
    13 = 132
    14 = 150
    15 = 4292
    18 = 467
    19 = 409
    20 = 985
    21 = 1549
    23 = 43
    25 = 45
    26 = 2004
    27 = 67
    28 = 62
    29 = 2
    30 = 136
    31 = 56
    32 = 204
    33 = 27
}

foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
